$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the client groups: append new client ids to the existing lists
# Row 10 -> RECURRENTE_>_500K    : add 20282
# Row 9  -> RECURRENTE_400K-500K : add 20205
$ws.Range("B10").Value = "11.20096.20228.20238.20256.20260.20303.20361.40124.50805.60125.60191.60217.60253.61012.90504.90602.90671.90509.20282"
$ws.Range("B9").Value = "955.993.998.1001.1006.1009.10424.20103.20125.20310.20384.40151.50623.60159.60162.60192.60225.70103.70113.10425.60126.40139.20205"

# Update the view: scroll so column B starts at the top-left and select B10
$ws.Range("B10").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
